# Update "想去人数" (want-to-go count) values in column F across sheets,
# reflecting a refreshed data snapshot (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 9380
$ws1.Range("F21").Value = 89
$ws1.Range("F36").Value = 190
$ws1.Range("F39").Value = 725

# --- Sheet: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F16").Value = 0
$ws2.Range("F22").Value = 649
$ws2.Range("F24").Value = 288
$ws2.Range("F25").Value = 288

# --- Sheet: 本地生活 (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 764
$ws3.Range("F7").Value = 2314
$ws3.Range("F8").Value = 3478

# --- Sheet: 全部类型 (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 9380
$ws4.Range("F8").Value = 3478
$ws4.Range("F36").Value = 288
$ws4.Range("F39").Value = 725
